$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set "Obrigatorio" (column E) to "S" for rows 2-9 and 11-15
$ws.Range("E2:E9").Value = "S"
$ws.Range("E11:E15").Value = "S"
